$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.929.44'
$ws.Range("E2").Value = '  -3.83%  '
$ws.Range("D3").Value = '1.637.31'
$ws.Range("E3").Value = '  -6.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9969'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.41%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4721'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2551'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06003'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07014'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.35%  '
$ws.Range("D11").Value = '1.638.11'
$ws.Range("E11").Value = '  -6.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6154'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.345'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '72.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.0000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9985'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '24.937.47'
$ws.Range("E18").Value = '  -3.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006576'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.12%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.846.38'
$ws.Range("E21").Value = '  -6.52%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.392'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.589'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.262'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '133.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.367'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '102.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.657'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.741'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07714'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.556'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9988'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04301'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.598'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9198'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5810'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.563'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01550'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9983'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8285'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.798'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3710'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.737'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.27%  '
$ws.Range("E46").Value = '  -3.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05213'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.076'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9993'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9976'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.64%  '
